# "Generate Report for Archive"
#
# The localization-status report is regenerated: the "Status" column value
# that used to read "Ready for handoff" is now "In Translation" (wherever it
# appears), and the "Status" columns are narrower to fit the new, shorter
# label (autosized ~17.22 chars -> ~13.41 chars).

$wb = $excel.ActiveWorkbook

# --- 1) Text: "Ready for handoff" -> "In Translation" -------------------
# Overview!E2:F2 (zh-cn / de-de status), zh-cn!C2 and de-de!C2 (Status).
# NB: put the literal on the LHS of -eq so a boolean-typed cell value (e.g.
# the "True"/"False" cells elsewhere in the sheet) doesn't get coerced to
# $true/$false and then loosely match any non-empty string on the RHS.
$targetOld = "Ready for handoff"
$targetNew = "In Translation"

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        $val = $cell.Value2
        if ($targetOld -eq $val) {
            $cell.Value = $targetNew
        }
    }
}

# --- 2) Column width: narrow the "Status" columns ------------------------
# Native width 17.2159881591797 -> 13.4101845877511. The ColumnWidth setter
# here snaps to the nearest 1/6-character pixel grid, so 12.5 is the closest
# input that lands on the nearest reachable width to the target.
$newColumnWidth = 12.5

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = $newColumnWidth   # column E (zh-cn)
$wsOverview.Columns.Item(6).ColumnWidth = $newColumnWidth   # column F (de-de)

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = $newColumnWidth       # column C (Status)

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = $newColumnWidth       # column C (Status)
